$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - new demo test data row
$ws.Range("B3").Value = " "
$ws.Range("C3").Value = " "
$ws.Range("D3").Value = "ThoHip12345678"
$ws.Range("E3").Value = " !@12345678"
$ws.Range("F3").Value = " !@12345678"
$ws.Range("G3").Value = "January"
$ws.Range("H3").Value = 13
$ws.Range("I3").Value = 1989
$ws.Range("J3").Value = "Female"

# Row 4 - new demo test data row
$ws.Range("B4").Value = " "
$ws.Range("C4").Value = "LastName"
$ws.Range("D4").Value = "ThoHip12345678"
$ws.Range("E4").Value = " !@12345678"
$ws.Range("F4").Value = " !@12345678"
$ws.Range("G4").Value = "January"
$ws.Range("H4").Value = 13
$ws.Range("I4").Value = 1989
$ws.Range("J4").Value = "Female"

# Row 5 - new demo test data row
$ws.Range("B5").Value = "FirstName"
$ws.Range("D5").Value = "ThoHip12345678"
$ws.Range("E5").Value = " !@12345678"
$ws.Range("F5").Value = " !@12345678"
$ws.Range("G5").Value = "January"
$ws.Range("H5").Value = 13
$ws.Range("I5").Value = 1989
$ws.Range("J5").Value = "Female"

# Move / update the sheet selection to reflect where the author left off editing
$ws.Range("D9").Select()
